# ---------------------------------------------------------------------------
# Apply the "updated results and plots" edit to all5.xlsx:
#   - add a new "wait-and-see" policy column (inserted as 3rd data column in
#     each of the two side-by-side tables)
#   - refresh the benchmark numbers for every existing column
#   - add matching relative-gap formulas/formatting for the new column
#   - the now-unused H:I filler cells below row 13 are cleared out
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Header row: shift old C/H ("RH2SSP-R") and D/I ("static2SSP-R")
#         one column to the right and insert "wait-and-see" in their place ----
$oldC1 = $ws.Range("C1").Value2
$oldD1 = $ws.Range("D1").Value2
$oldH1 = $ws.Range("H1").Value2
$oldI1 = $ws.Range("I1").Value2

$ws.Range("E1").Value = $oldD1
$ws.Range("D1").Value = $oldC1
$ws.Range("C1").Value = "wait-and-see"

$ws.Range("J1").Value = $oldI1
$ws.Range("I1").Value = $oldH1
$ws.Range("H1").Value = "wait-and-see"

# ---- 2. Refreshed benchmark data, rows 2-11, columns A-E ----
# A=CV-R  B=FA-MSP-R  C=wait-and-see  D=RH2SSP-R  E=static2SSP-R  F=gap (unchanged)
$data = @(
  @(6301.9963310000003, 6377.7896479999999, 6395.5240489999996, 8511.721254,        22819.79292),
  @(6663.3840319999999, 7159.1757530000004, 7528.6315869999999, 8716.0863879999997, 22938.59158),
  @(6885.9000150000002, 7712.2349869999998, 8661.7391239999997, 8920.4515210000009, 23057.390240000001),
  @(7090.8903280000004, 8183.7728230000002, 9794.8466609999996, 9124.8166550000005, 23176.188890000001),
  @(7295.2690009999997, 8598.2645190000003, 10927.9542,         9329.1817890000002, 23294.987550000002),
  @(7499.6341350000002, 8988.9785859999993, 12061.061739999999, 9533.5469229999999, 23413.786209999998),
  @(7703.9992689999999, 9366.2381389999991, 13194.16927,        9737.9120569999995, 23532.584869999999),
  @(7908.3644020000002, 9670.3725030000005, 14327.276809999999, 9942.2771900000007, 23651.383529999999),
  @(8112.7295359999998, 9950.084476,        15460.38435,        10146.642320000001, 23770.18219),
  @(8317.0946700000004, 10231.2765,         16593.491890000001, 10351.007460000001, 23888.98085)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]   # A
  $ws.Cells.Item($r, 2).Value = $row[1]   # B
  $ws.Cells.Item($r, 3).Value = $row[2]   # C
  $ws.Cells.Item($r, 4).Value = $row[3]   # D
  $ws.Cells.Item($r, 5).Value = $row[4]   # E
}

# Column E previously held unused, centre-aligned "0.00" placeholder
# formatting; the new static2SSP-R figures use plain general formatting.
$ws.Range("E2:E11").ClearFormats()
for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 5).Value = $data[$i][4]
}

# ---- 3. Relative-gap formulas, rows 2-11, columns G-J ----
# G (vs B) already existed and keeps its formula - just needs the new
# number format below. H/I keep referencing C/D (whose meaning shifted to
# wait-and-see / RH2SSP-R) and J is brand new, referencing E.
for ($r = 2; $r -le 11; $r++) {
  $ws.Cells.Item($r, 8).Formula = "=(C$r-`$A$r)/`$A$r"   # H
  $ws.Cells.Item($r, 9).Formula = "=(D$r-`$A$r)/`$A$r"   # I
  $ws.Cells.Item($r, 10).Formula = "=(E$r-`$A$r)/`$A$r"  # J
}

# ---- 4. Number formatting for the gap columns (one-decimal percent) ----
$ws.Range("G2:J11").NumberFormat = "0.0%"

# ---- 5. The filler rows below the table no longer carry H/I formatting
#         past row 13 (only F/G remain styled placeholders) ----
$ws.Range("H14:I52").Clear()

$wb.Save()
